$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 11 row was previously blank (aside from the week number and the
# H11 shared-formula total). Fill in the date and the five contributors'
# percentages, matching the pattern used by the other completed weeks.
$ws.Range("B11").Value = 43920

# Copy the date formatting from the row above (B10) onto B11 so it picks
# up the existing date number format style instead of minting a new one.
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C11").Value = 0.2
$ws.Range("D11").Value = 0.2
$ws.Range("E11").Value = 0.2
$ws.Range("F11").Value = 0.2
$ws.Range("G11").Value = 0.2

# Move the active selection to C14, matching where the user clicked next.
$ws.Range("C14").Select()
